$wb = $excel.ActiveWorkbook

# Update the image path text: images now live under an "images/" subfolder
$mainWs = $wb.Worksheets.Item("Main")
$mainWs.Range("G2").Value = "/recordm/localresource/dash/assets/images/improve.jpg"

# The ">Text" sheet was the active/selected one; move the active selection
# to the "Main" sheet (cell G2, the one that was just edited) instead.
$mainWs.Activate() | Out-Null
$mainWs.Range("G2").Select() | Out-Null
